# Commit: "Sun, Apr 26, 2020  7:06:26 PM"
#
# Two changes are applied to the deck:
#
#  1. The table on slide 16 gets a different built-in table style
#     (from the deck's custom "Table_0" style to PowerPoint's gallery
#     style {8019AB9C-3FCF-43EC-AD14-B6139F48F853}).
#
#  2. The presentation's applied theme switches from the custom
#     "Integral" design back to the stock "Office Theme" palette
#     (Design tab -> Themes -> Office Theme). We reproduce this by
#     writing the twelve standard Office Theme RGB values into the
#     Slide Master's theme color scheme, which is the live theme for
#     every slide in the deck.

$p = $ppt.ActivePresentation

# --- 1. Table style on slide 16 --------------------------------------
$slide = $p.Slides.Item(16)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{8019AB9C-3FCF-43EC-AD14-B6139F48F853}")
    }
}

# --- 2. Switch theme colors: "Integral" -> "Office Theme" ------------
function Set-ThemeColor($scheme, $index, $hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    $rgbValue = $r + ($g * 256) + ($b * 65536)
    $scheme.Colors($index).RGB = $rgbValue
}

$colorScheme = $p.SlideMaster.ColorScheme

# Office Theme palette, in dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink order
Set-ThemeColor $colorScheme 1  "000000"
Set-ThemeColor $colorScheme 2  "FFFFFF"
Set-ThemeColor $colorScheme 3  "44546A"
Set-ThemeColor $colorScheme 4  "E7E6E6"
Set-ThemeColor $colorScheme 5  "5B9BD5"
Set-ThemeColor $colorScheme 6  "ED7D31"
Set-ThemeColor $colorScheme 7  "A5A5A5"
Set-ThemeColor $colorScheme 8  "FFC000"
Set-ThemeColor $colorScheme 9  "4472C4"
Set-ThemeColor $colorScheme 10 "70AD47"
Set-ThemeColor $colorScheme 11 "0563C1"
Set-ThemeColor $colorScheme 12 "954F72"
